$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (rows 2:317) holds one price record per row, sorted by date
# within each variety/quality group. A new record needs to be inserted right
# before the current row 244, which pushes every following row down by one.
$ws.Rows("244:244").Insert()

# Seed the freshly inserted row 244 with the same record that now lives in
# row 245 (i.e. the record that used to be row 244 before the insert), since
# the new entry shares almost all of its attributes with it.
$ws.Range("A245:T245").Copy()
$ws.Range("A244:T244").PasteSpecial()

# Now adjust the handful of fields that differ for the new record: the date,
# the quality grade and the traded volume.
$ws.Range("D244").Value = 44841
$ws.Range("L244").Value = "Primera"
$ws.Range("M244").Value = 100
